$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "memory leak"
$ws.Range("D11").Value = "valgrind"
$ws.Range("B13").Value = "arm assemble"

$ws.Range("B14").Select()
